# Applies the "adjustments for next simulation in hyperparam tuning,
# added NN architecture" edit to the Optimizer_Testing_NN_outcome workbook.
#
# - Fills in the previously-blank Time/RMSE result columns (C:R) for the
#   per-repetition rows on both the Single_Layer and Double_Layer sheets
#   (formulas in rows 13-16 / 28-31 recalc automatically).
# - Adjusts a couple of cell highlight styles that moved along with the
#   new numbers.
# - Updates the saved cursor/selection on both sheet views.

$wb = $excel.ActiveWorkbook

function Set-RowValues {
    param(
        $ws,
        [int]$row,
        [int]$startCol,
        [double[]]$vals
    )
    $n = $vals.Length
    $arr = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0, $i] = $vals[$i]
    }
    $endCol = $startCol + $n - 1
    $rng = $ws.Range($ws.Cells.Item($row, $startCol), $ws.Cells.Item($row, $endCol))
    $rng.Value = $arr
}

$ws1 = $wb.Worksheets.Item("Single_Layer")
$ws2 = $wb.Worksheets.Item("Double_Layer")

# ---------------------------------------------------------------------
# Single_Layer: rows 8-11, columns C:R
# (NB: call the helper positionally -- named-parameter binding isn't
# reliable for script functions in this host)
# ---------------------------------------------------------------------
Set-RowValues $ws1 8  3 @(8.09, 1.13, 6, 1.13, 4.63, 1.12, 4.39, 1.12, 13.97, 1.16, 10.86, 1.14, 12.16, 1.12, 10.18, 1.12)
Set-RowValues $ws1 9  3 @(5.07, 1.13, 4.6, 1.13, 4.68, 1.13, 3.56, 1.14, 14.46, 1.14, 10.23, 1.14, 11.66, 1.14, 8.37, 1.17)
Set-RowValues $ws1 10 3 @(5.09, 1.14, 5.31, 1.14, 4.15, 1.16, 3.62, 1.17, 8.68, 1.14, 6.5, 1.13, 7.28, 1.13, 6.58, 1.13)
Set-RowValues $ws1 11 3 @(20.64, 1.13, 16.83, 1.12, 13.99, 1.11, 13.3, 1.11, 30.47, 1.2, 20.66, 1.21, 19.45, 1.21, 13.15, 1.21)

# I11:J11 pick up the "min" highlight fill (style index 3 in the original file)
$ws1.Range("I11:J11").Interior.Color = 65535

# ---------------------------------------------------------------------
# Double_Layer: rows 7-11 and 22-26, columns C:R
# ---------------------------------------------------------------------
Set-RowValues $ws2 7  3 @(4.52, 1.14, 4.56, 1.16, 4.36, 1.16, 3.8, 1.19, 6.14, 1.14, 5.42, 1.15, 5.63, 1.15, 6.38, 1.15)
Set-RowValues $ws2 8  3 @(7.72, 1.14, 5.75, 1.16, 7.37, 1.14, 5.31, 1.15, 9.02, 1.16, 9.2, 1.13, 8.31, 1.14, 8.72, 1.13)
Set-RowValues $ws2 9  3 @(4.92, 1.15, 4.4, 1.15, 4.5, 1.15, 4.13, 1.16, 6.97, 1.15, 6.26, 1.14, 6.28, 1.14, 7.31, 1.13)
Set-RowValues $ws2 10 3 @(5.29, 1.16, 4.93, 1.17, 4.72, 1.17, 4.78, 1.19, 7.36, 1.15, 6.3, 1.16, 6.16, 1.15, 5.65, 1.15)
Set-RowValues $ws2 11 3 @(16.8, 1.15, 16, 1.14, 19.91, 1.14, 11.95, 1.15, 21.48, 1.21, 15.73, 1.2, 23.7, 1.17, 17.39, 1.18)

Set-RowValues $ws2 22 3 @(4.99, 1.15, 4.06, 1.15, 6.01, 1.14, 4.25, 1.15, 6.43, 1.15, 6.56, 1.14, 6.81, 1.15, 6.01, 1.15)
Set-RowValues $ws2 23 3 @(8.16, 1.13, 8.18, 1.14, 6.51, 1.14, 6.31, 1.13, 14.87, 1.15, 12.08, 1.13, 14.1, 1.16, 12.6, 1.14)
Set-RowValues $ws2 24 3 @(5.49, 1.15, 4.96, 1.13, 5.6, 1.13, 4.78, 1.13, 16.77, 1.15, 12.3, 1.14, 14.72, 1.15, 11.35, 1.16)
Set-RowValues $ws2 25 3 @(5, 1.14, 5.44, 1.15, 5.25, 1.15, 4.93, 1.15, 10.6, 1.15, 7.83, 1.14, 9.27, 1.15, 8.58, 1.14)
Set-RowValues $ws2 26 3 @(25.82, 1.2, 23.46, 1.13, 25.23, 1.14, 29.06, 1.13, 22.11, 1.28, 31.82, 1.21, 27.6, 1.23, 14.18, 1.31)

# Q8:R8 loses the highlight fill it previously had (back to the plain bordered style)
$src = $ws2.Range("P8")
$src.Copy()
$ws2.Range("Q8:R8").PasteSpecial(-4122)

# C23:D23 drops its border/style entirely (reverts to the default, unstyled cell)
$src2 = $ws2.Range("D13")
$src2.Copy()
$ws2.Range("C23:D23").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Selection / cursor updates recorded in the sheet views
# ---------------------------------------------------------------------
$ws1.Range("W10").Select()

$ws2.Application.Goto($ws2.Range("A7"), $false)
$ws2.Range("J24").Select()
